$d = $word.ActiveDocument

function Get-ParaIndexContaining {
    param($doc, $needle)
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

function Merge-ParagraphRuns {
    # Collapses every run in the paragraph at $paraIndex into a single run
    # that carries the combined text, preserving italic/bold character
    # formatting if the original (first) run had it.
    param($doc, $paraIndex)
    $p = $doc.Paragraphs.Item($paraIndex)
    $full = $p.Range
    $textRange = $doc.Range($full.Start, $full.End - 1)
    $combinedText = $textRange.Text
    $isItalic = $textRange.Italic
    $isBold = $textRange.Bold
    $textRange.Delete()
    $p.Range.InsertBefore($combinedText)
    $newRange = $doc.Range($p.Range.Start, $p.Range.End - 1)
    if ($isItalic -eq -1) {
        $newRange.Font.Italic = -1
    }
    if ($isBold -eq -1) {
        $newRange.Font.Bold = -1
    }
}

# 1) "Católico/Protestante: " + "Religion;" -> single italic run.
$idxCatolico = Get-ParaIndexContaining $d "Católico/Protestante"
Merge-ParagraphRuns $d $idxCatolico

# 2) Insert the new bold "INTERFACE" block right after the
#    "Ações (TODO):" paragraph, before the "Unique name" bullet list.
#    (Search on "TODO" only -- the source text has a non-breaking space
#    between "Ações" and "(TODO):" that a literal match would miss.)
$idxAcoes = Get-ParaIndexContaining $d "TODO"
$acoesPara = $d.Paragraphs.Item($idxAcoes)
$acoesPara.Range.InsertParagraphAfter()

$interfaceLines = @(
    "INTERFACE",
    "- desaparecer a tela inicial de adicionar jogadores",
    "- mostrar as cartas de cada jogador na mesa",
    "- indicar o jogador atual",
    "- mostrar ações possíveis para o jogador atual"
)

$curIdx = $idxAcoes + 1
for ($i = 0; $i -lt $interfaceLines.Count; $i++) {
    $curPara = $d.Paragraphs.Item($curIdx)
    $curPara.Range.InsertAfter($interfaceLines[$i])
    if ($i -lt ($interfaceLines.Count - 1)) {
        $curPara.Range.InsertParagraphAfter()
    }
    $curIdx = $curIdx + 1
}
# Trailing empty bold paragraph after the INTERFACE block.
$d.Paragraphs.Item($curIdx - 1).Range.InsertParagraphAfter()

# 3) "Unique name – dois " + "jogadores não podem ter o mesmo nome;" -> single run.
$idxUnique = Get-ParaIndexContaining $d "Unique name"
Merge-ParagraphRuns $d $idxUnique

# 4) "V" + "erificar regulamentação;" -> single run.
$idxVerificar = Get-ParaIndexContaining $d "erificar regulament"
Merge-ParagraphRuns $d $idxVerificar

Write-Output "done"
